$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, pushing the existing row 47 (and all
# rows below it) down by one. This mirrors the diff: a new "Poroto
# granado" price record is added after row 46, and every subsequent
# record shifts down one row (old row 86 becomes new row 87).
$ws.Rows("47:47").Insert()

# Populate the newly inserted row 47 with the new record's data. The
# fixed/common columns (A,B,C,E,F,G,H,I,N,O,Q,R) match every other row
# in this sheet for "Macroferia Regional de Talca" / "Poroto granado".
$ws.Range("A47").Value = 5
$ws.Range("B47").Value = "Macroferia Regional de Talca"
$ws.Range("C47").Value = "Maule"
$ws.Range("D47").Value = 44566
$ws.Range("E47").Value = 7
$ws.Range("F47").Value = 100112030
$ws.Range("G47").Value = "Poroto granado"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 300
$ws.Range("K47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("M47").Value = 30000
$ws.Range("N47").Value = "`$/saco 25 kilos"
$ws.Range("O47").Value = "Región del Maule"
$ws.Range("P47").Value = 1200
$ws.Range("Q47").Value = 25
$ws.Range("R47").Value = "Hortaliza"
